# Loan RBI, Variable Instalments
# Insert a new column before column N ("Late") on the "Repayment schedule"
# sheet, then leave the "Repayment schedule" sheet active/selected (as the
# Transactions sheet used to be), matching what the author did in Excel.

$wb = $excel.ActiveWorkbook

$repayment = $wb.Worksheets.Item("Repayment schedule")
$transactions = $wb.Worksheets.Item("Transactions")

# Insert a new column before column N (shifts Late/heading/Outstanding right)
$repayment.Columns("N").Insert()

# Selections: Transactions sheet loses its "selected"/tabSelected state and
# moves the active cell/selection to column D; Repayment schedule becomes the
# active sheet with a new selection.
$transactions.Range("D1:D1048576").Select()

$repayment.Activate()
$repayment.Range("J17").Select()
